$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "1" to "Khashuri"
$ws.Name = "Khashuri"

# Row 6 (Urban) - some values became confidential ("...")
$row6Dots = @("C","E","F","I","K","N","O")
foreach ($col in $row6Dots) {
    $ws.Range($col + "6").Value = "..."
}

# Row 7 (Rural) - some values became confidential ("...")
$row7Dots = @("C","E","F","H","I","J","K","N","O")
foreach ($col in $row7Dots) {
    $ws.Range($col + "7").Value = "..."
}

# Row 8 was already empty; the note that was on row 9 moves up to row 8
$ws.Rows("8").Delete()
